$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear the "Detail" notes that no longer apply
$ws.Range("I33").Value = $null
$ws.Range("I36").Value = $null

# Fill in In-charge Roll Number / In-charge Full Name for "Posts List" and "Post Details"
# rows, matching the values already used for "Blogs List"/"Blog Details" (HE150340 / Phung Quang Thong)
$ws.Range("E36").Value = "HE150340"
$ws.Range("F36").Value = "Phung Quang Thong(BE-50, FE-50)"
$ws.Range("E37").Value = "HE150340"
$ws.Range("F37").Value = "Phung Quang Thong(BE-50, FE-50)"

# Fill in In-charge Roll Number / In-charge Full Name for "Feedbacks List" and "Feedback Details"
$ws.Range("E40").Value = "HE150411"
$ws.Range("F40").Value = "Tran Tat Dat(BE-50, FE-50)"
$ws.Range("E41").Value = "HE150411"
$ws.Range("F41").Value = "Tran Tat Dat(BE-50, FE-50)"

# Update the view: scroll to show column D and select I33
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 4
$ws.Range("I33").Select()
